# Make styling consistent across metadata for each lesson
#
# The "we want to hear from you" paragraph (ending in the
# petascale@shodor.org hyperlink) was immediately followed by its own
# empty paragraph that only existed to carry a page break run. That
# paragraph break is removed so the page-break run joins the hyperlink
# paragraph, and the document's "_GoBack" bookmark (which Word had
# stranded at the very end of the document) is relocated to sit right at
# that newly-created join point instead.

$d = $word.ActiveDocument

# --- Step 1: drop the existing "_GoBack" bookmark whichever paragraph it
#     currently lives in (at the very end of the document). ---
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- Step 2: find the paragraph that ends with the petascale@shodor.org
#     hyperlink. ---
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*petascale@shodor.org*") {
        $targetIndex = $i
    }
}

$para = $d.Paragraphs($targetIndex)
$joinPoint = $para.Range.End - 1

# --- Step 3: merge that paragraph with the following (page-break-only)
#     paragraph by deleting the paragraph mark between them. ---
$markRange = $d.Range($joinPoint, $joinPoint + 1)
$markRange.Delete()

# --- Step 4: re-create the "_GoBack" bookmark right at the junction. ---
$bmRange = $d.Range($joinPoint, $joinPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
